$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A=6) -> name becomes "line7"
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11

# Row 9 (A=7) -> name becomes "line8"
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16

# Row 10 (A=8) -> name becomes "extr1"
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (A=9) -> name becomes "extr2"
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9

# Row 12 (A=10) -> name becomes "extr3"
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10

# Row 13 (A=11) -> name becomes "extr4"
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

# Row 14 (A=12) -> name becomes "extr5"
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

# Row 15 (A=13) -> name becomes "extr6"
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# New row 16 (A=14) -> name "extr7"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

# New row 17 (A=15) -> name "extr8"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

# Copy column-A formatting (bold, centered, thin border) from the last
# existing data row onto the two newly appended rows so they look
# consistent with the rest of the table, then restore the numeric value
# (PasteSpecial formats-only keeps the existing value untouched, but we
# set it again defensively).
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(17, 1).Value = 15
